$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing song title in B3
$ws.Range("B3").Value = "I Have Decided to Follow Jesus [Kuv Twv Xais Tag Lawm] (#95)"

# Adjust column B width to fit new content (target stored width 56.5703125;
# ColumnWidth is offset by ~5/6 before being stored, so compensate here)
$ws.Columns.Item(2).ColumnWidth = 55.666666666666664

# Update the active selection to B4
$ws.Range("B4").Select()
